# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "last updated" timestamp text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 20 de Agosto de 2020 a las 08:05"

# --- Update numeric stats for several countries ---

# Row 6: India
$ws.Range("B6").Value = 2837749
$ws.Range("C6").Value = 1927
$ws.Range("D6").Value = 2097331
$ws.Range("E6").Value = 686414
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 54004

# Row 18: Pakistan
$ws.Range("B18").Value = 290958
$ws.Range("C18").Value = 513
$ws.Range("D18").Value = 272804
$ws.Range("E18").Value = 11945
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = 6209

# Row 32: Israel
$ws.Range("B32").Value = 98443
$ws.Range("C32").Value = 474
$ws.Range("D32").Value = 73841
$ws.Range("E32").Value = 23821

# Row 56: Kirguistan
$ws.Range("B56").Value = 42507
$ws.Range("C56").Value = 182
$ws.Range("D56").Value = 35486
$ws.Range("E56").Value = 5523

# Row 61: Uzbekistan
$ws.Range("B61").Value = 37366
$ws.Range("C61").Value = 254
$ws.Range("E61").Value = 4172
$ws.Range("G61").Value = 2
$ws.Range("H61").Value = 250

# Row 71: Australia
$ws.Range("D71").Value = 15249
$ws.Range("E71").Value = 8524

# Row 73: El Salvador
$ws.Range("D73").Value = 11428
$ws.Range("E73").Value = 11649
$ws.Range("G73").Value = 7
$ws.Range("H73").Value = 640

# Row 161: Reunion
$ws.Range("B161").Value = 945
$ws.Range("E161").Value = 283

# --- Swap "Islas Malvinas" (row 213) and "Montserrat" (row 214) ---
# The shared-string order for these two countries was swapped upstream,
# which (combined with the updated stat values) results in the two rows
# exchanging their full contents.
$row213 = @()
for ($c = 1; $c -le 8; $c++) {
    $row213 += , $ws.Cells.Item(213, $c).Value()
}
$row214 = @()
for ($c = 1; $c -le 8; $c++) {
    $row214 += , $ws.Cells.Item(214, $c).Value()
}
for ($c = 1; $c -le 8; $c++) {
    $ws.Cells.Item(213, $c).Value = $row214[$c - 1]
    $ws.Cells.Item(214, $c).Value = $row213[$c - 1]
}

Write-Output "Edit complete"
